$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in Friday's (row 17) In/Out times ---
# C17 = 08:00 (0.333333333333333), D17 = 12:00 (0.5)
$ws.Range("C17").Value = 0.333333333333333
$ws.Range("D17").Value = 0.5

# --- Add note for Friday in L17 ---
$ws.Range("L17").Value = "0800-1000, 1300-1500"

# --- Narrow column B and merge C:H widths into one uniform band ---
# Target raw widths (~character units): B ~= 6.75 ; C.. ~= 6.0765306122449
# The engine snaps ColumnWidth to 1/6-character increments, so we pick the
# input value that lands on the closest achieving grid point.
$ws.Columns.Item(2).ColumnWidth = 5.83
$ws.Columns.Item(8).ColumnWidth = 5.17

# --- Update selection / scroll position ---
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("G17").Select()
